# issue #5: add legislator_id, name, date into dataframe
#
# The workbook's "股票" (stocks) sheet gains three new trailing columns:
#   H: date             -> "2012-02-13" for every data row
#   I: legislator_name  -> "陳明文" for every data row
#   J: legislator_id    -> 828 for every data row
#
# Header formatting should mirror the existing header cells (bold, centred,
# bordered) and data-row formatting should mirror the existing data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- header row -----------------------------------------------------------
# Copy the formatting of the last existing header cell (G1) onto the three
# new header cells, then set their text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# --- data rows --------------------------------------------------------------
$lastRow = 8

for ($r = 2; $r -le $lastRow; $r++) {
    # Mirror the formatting used by the rest of the row (column G) onto the
    # three new cells for this row.
    $ws.Range("G$r").Copy()
    $ws.Range("H$r").PasteSpecial(-4122)
    $ws.Range("I$r").PasteSpecial(-4122)
    $ws.Range("J$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 8).Value = "2012-02-13"
    $ws.Cells.Item($r, 9).Value = "陳明文"
    $ws.Cells.Item($r, 10).Value = 828
}
